$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A338").Value = 337
$ws.Range("B338").Value = "Thursday, Jan 12"
$ws.Range("C338").Value = "8:10 AM"
$ws.Range("D338").Value = "LO3903"
$ws.Range("E338").Value = "Warsaw"
$ws.Range("F338").Value = "(WAW)"
$ws.Range("G338").Value = "LOT "
$ws.Range("H338").Value = "E170"
$ws.Range("I338").Value = "(SP-LDF)"
$ws.Range("J338").Value = "8:00 AM"
$ws.Range("L338").Value = "0 hours, -10 minutes"

$ws.Range("A339").Value = 338
$ws.Range("B339").Value = "Thursday, Jan 12"
$ws.Range("C339").Value = "8:45 AM"
$ws.Range("D339").Value = "FR3510"
$ws.Range("E339").Value = "Milan"
$ws.Range("F339").Value = "(BGY)"
$ws.Range("G339").Value = "Ryanair "
$ws.Range("H339").Value = "B38M"
$ws.Range("I339").Value = "(9H-VUK)"
$ws.Range("J339").Value = "8:31 AM"
$ws.Range("L339").Value = "0 hours, -14 minutes"

$ws.Range("A340").Value = 339
$ws.Range("B340").Value = "Thursday, Jan 12"
$ws.Range("C340").Value = "9:45 AM"
$ws.Range("D340").Value = "FR6293"
$ws.Range("E340").Value = "London"
$ws.Range("F340").Value = "(STN)"
$ws.Range("G340").Value = "Ryanair "
$ws.Range("H340").Value = "B38M"
$ws.Range("I340").Value = "(EI-HEW)"
$ws.Range("J340").Value = "9:27 AM"
$ws.Range("L340").Value = "0 hours, -18 minutes"

$ws.Range("A341").Value = 340
$ws.Range("B341").Value = "Thursday, Jan 12"
$ws.Range("C341").Value = "9:55 AM"
$ws.Range("D341").Value = "RK6249"
$ws.Range("E341").Value = "Manchester"
$ws.Range("F341").Value = "(MAN)"
$ws.Range("G341").Value = "Ryanair "
$ws.Range("H341").Value = "B738"
$ws.Range("I341").Value = "(G-RUKD)"
$ws.Range("J341").Value = "9:32 AM"
$ws.Range("L341").Value = "0 hours, -23 minutes"

$ws.Range("A342").Value = 341
$ws.Range("B342").Value = "Thursday, Jan 12"
$ws.Range("C342").Value = "10:15 AM"
$ws.Range("D342").Value = "FR5118"
$ws.Range("E342").Value = "Dublin"
$ws.Range("F342").Value = "(DUB)"
$ws.Range("G342").Value = "Ryanair "
$ws.Range("H342").Value = "B38M"
$ws.Range("I342").Value = "(EI-HES)"
$ws.Range("J342").Value = "9:57 AM"
$ws.Range("L342").Value = "0 hours, -18 minutes"

$ws.Range("A343").Value = 342
$ws.Range("B343").Value = "Thursday, Jan 12"
$ws.Range("C343").Value = "11:00 AM"
$ws.Range("D343").Value = "FR6211"
$ws.Range("E343").Value = "Paris"
$ws.Range("F343").Value = "(BVA)"
$ws.Range("G343").Value = "Buzz "
$ws.Range("H343").Value = "B38M"
$ws.Range("I343").Value = "(SP-RZA)"
$ws.Range("J343").Value = "10:57 AM"
$ws.Range("L343").Value = "0 hours, -3 minutes"

$ws.Range("A344").Value = 343
$ws.Range("B344").Value = "Thursday, Jan 12"
$ws.Range("C344").Value = "11:30 AM"
$ws.Range("D344").Value = "LO3907"
$ws.Range("E344").Value = "Warsaw"
$ws.Range("F344").Value = "(WAW)"
$ws.Range("G344").Value = "LOT "
$ws.Range("H344").Value = "B38M"
$ws.Range("I344").Value = "(SP-LVA)"
$ws.Range("J344").Value = "11:18 AM"
$ws.Range("L344").Value = "0 hours, -12 minutes"

$ws.Range("A345").Value = 344
$ws.Range("B345").Value = "Thursday, Jan 12"
$ws.Range("C345").Value = "11:35 AM"
$ws.Range("D345").Value = "FR6257"
$ws.Range("E345").Value = "Stockholm"
$ws.Range("F345").Value = "(ARN)"
$ws.Range("G345").Value = "Buzz "
$ws.Range("H345").Value = "B38M"
$ws.Range("I345").Value = "(SP-RZD)"
$ws.Range("J345").Value = "11:25 AM"
$ws.Range("L345").Value = "0 hours, -10 minutes"

$ws.Range("A346").Value = 345
$ws.Range("B346").Value = "Thursday, Jan 12"
$ws.Range("C346").Value = "11:50 AM"
$ws.Range("D346").Value = "W65060"
$ws.Range("E346").Value = "Catania"
$ws.Range("F346").Value = "(CTA)"
$ws.Range("G346").Value = "Wizz Air "
$ws.Range("H346").Value = "A321"
$ws.Range("I346").Value = "(HA-LXO)"
$ws.Range("J346").Value = "11:29 AM"
$ws.Range("L346").Value = "0 hours, -21 minutes"

$ws.Range("A347").Value = 346
$ws.Range("B347").Value = "Thursday, Jan 12"
$ws.Range("C347").Value = "11:55 AM"
$ws.Range("D347").Value = "D84901"
$ws.Range("E347").Value = "Stockholm"
$ws.Range("F347").Value = "(ARN)"
$ws.Range("G347").Value = "Norwegian "
$ws.Range("H347").Value = "B738"
$ws.Range("I347").Value = "(SE-RPE)"
$ws.Range("J347").Value = "11:33 AM"
$ws.Range("L347").Value = "0 hours, -22 minutes"

$ws.Range("A348").Value = 347
$ws.Range("B348").Value = "Thursday, Jan 12"
$ws.Range("C348").Value = "11:55 AM"
$ws.Range("D348").Value = "LH1366"
$ws.Range("E348").Value = "Frankfurt"
$ws.Range("F348").Value = "(FRA)"
$ws.Range("G348").Value = "Lufthansa "
$ws.Range("H348").Value = "A21N"
$ws.Range("I348").Value = "(D-AIEP)"
$ws.Range("J348").Value = "12:24 PM"
$ws.Range("L348").Value = "0 hours, 29 minutes"

$ws.Range("A349").Value = 348
$ws.Range("B349").Value = "Thursday, Jan 12"
$ws.Range("C349").Value = "12:05 PM"
$ws.Range("D349").Value = "DY1040"
$ws.Range("E349").Value = "Oslo"
$ws.Range("F349").Value = "(OSL)"
$ws.Range("G349").Value = "Norwegian "
$ws.Range("H349").Value = "B738"
$ws.Range("I349").Value = "(LN-DYJ)"
$ws.Range("J349").Value = "11:53 AM"
$ws.Range("L349").Value = "0 hours, -12 minutes"

$ws.Range("A350").Value = 349
$ws.Range("B350").Value = "Thursday, Jan 12"
$ws.Range("C350").Value = "12:25 PM"
$ws.Range("D350").Value = "FR1813"
$ws.Range("E350").Value = "London"
$ws.Range("F350").Value = "(LTN)"
$ws.Range("G350").Value = "Ryanair "
$ws.Range("H350").Value = "B738"
$ws.Range("I350").Value = "(SP-RKC)"
$ws.Range("J350").Value = "12:30 PM"
$ws.Range("L350").Value = "0 hours, 5 minutes"

$ws.Range("A351").Value = 350
$ws.Range("B351").Value = "Thursday, Jan 12"
$ws.Range("C351").Value = "12:35 PM"
$ws.Range("D351").Value = "U28511"
$ws.Range("E351").Value = "London"
$ws.Range("F351").Value = "(LGW)"
$ws.Range("G351").Value = "easyJet "
$ws.Range("H351").Value = "A319"
$ws.Range("I351").Value = "(G-EZBW)"
$ws.Range("J351").Value = "12:09 PM"
$ws.Range("L351").Value = "0 hours, -26 minutes"
